$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, derived from the authoritative commit diff.
# Covers: (1) item-row quantity (F) and extended-value (G = Rate * Qty) corrections,
# (2) company "Sub Total:" (B) rollups recomputed as SUM(G) over each group,
# (3) a handful of adjacent data-row pairs whose Code/MRP/Qty/Value (B/E/F/G) were
#     transposed between the two rows, and (4) the workbook grand totals (B793/B794).
$updates = [ordered]@{
    "F9" = 208
    "G9" = 6150.56
    "B10" = 46121.3
    "F83" = 82
    "G83" = 5476.78
    "F85" = 25
    "G85" = 568.25
    "F93" = 75
    "G93" = 5923.5
    "F97" = 385
    "G97" = 17994.9
    "F102" = 33
    "G102" = 3221.46
    "F109" = 165
    "G109" = 20702.55
    "B114" = 320411.49
    "F192" = 1
    "G192" = 65.48999999999999
    "F193" = 350
    "G193" = 22680
    "F197" = 116
    "G197" = 5392.84
    "F198" = 63
    "G198" = 4733.82
    "B200" = 60779
    "F219" = 17
    "G219" = 11424.68
    "B222" = 73918.14
    "F229" = 68
    "G229" = 2180.76
    "F247" = 138
    "G247" = 4595.4
    "B274" = 121898.93
    "F281" = 36
    "G281" = 5792.4
    "B292" = 55373
    "E292" = 163.62
    "F292" = -94
    "G292" = -13562.32
    "B293" = 63520
    "E293" = 153.4
    "F293" = 83
    "G293" = 11975.24
    "B294" = 57802
    "E294" = 162.71
    "F294" = -79
    "G294" = -11334.92
    "F295" = 11
    "G295" = 1578.28
    "B296" = 63531
    "E296" = 152.53
    "F296" = 80
    "G296" = 11478.4
    "F313" = 40
    "G313" = 4854.4
    "F319" = 11
    "G319" = 131.34
    "F322" = 222
    "G322" = 23068.02
    "F323" = 55
    "G323" = 6524.65
    "F325" = 195
    "G325" = 26874.9
    "B339" = 416337.8
    "F344" = 2
    "G344" = 416.6
    "B346" = 15760.32
    "F355" = 133
    "G355" = 4287.92
    "B361" = 16059.29
    "F363" = 57
    "G363" = 1211.25
    "F370" = 61
    "G370" = 3275.7
    "F376" = 83
    "G376" = 581
    "F379" = 98
    "G379" = 9113.02
    "F382" = 182
    "G382" = 7827.82
    "F389" = 6
    "G389" = 335.16
    "B395" = 281782.18
    "F404" = 0
    "G404" = 0
    "B411" = 4208.88
    "F424" = 90
    "G424" = 2676.6
    "F426" = 245
    "G426" = 23667
    "F427" = 128
    "G427" = 4762.88
    "B430" = 64607.76
    "B465" = 65069
    "E465" = 14.3
    "F465" = 23
    "G465" = 309.35
    "B466" = 53757
    "E466" = 16.08
    "F466" = -159
    "G466" = -2138.55
    "F468" = 146
    "G468" = 1919.9
    "F470" = 205
    "G470" = 2626.05
    "F477" = 156
    "G477" = 3077.88
    "B490" = 65067
    "E490" = 15.65
    "F490" = 271
    "G490" = 3991.83
    "B491" = 53595
    "E491" = 17.61
    "F491" = -335
    "G491" = -4934.55
    "B492" = 6949.68
    "F506" = 9
    "G506" = 626.85
    "B508" = 25681.97
    "F511" = 113
    "G511" = 2740.25
    "B528" = 25067.81
    "F547" = 174
    "G547" = 3453.9
    "B552" = 10466.92
    "B596" = 64830
    "E596" = 34.9
    "F596" = 114
    "G596" = 3742.62
    "B597" = 60022
    "E597" = 37.22
    "F597" = -113
    "G597" = -3709.79
    "F650" = 402
    "G650" = 32312.76
    "B651" = 42459.45
    "F698" = 107
    "G698" = 8726.92
    "F699" = 57
    "G699" = 2728.02
    "F701" = 250
    "G701" = 35782.5
    "F702" = 96
    "G702" = 7829.76
    "F704" = 38
    "G704" = 5057.8
    "F707" = 193
    "G707" = 4191.96
    "F713" = 576
    "G713" = 77765.75999999999
    "F714" = 60
    "G714" = 2245.2
    "F715" = 506
    "G715" = 61079.26
    "B716" = 245794.38
    "F718" = 47
    "G718" = 11725.56
    "F722" = 14
    "G722" = 1523.34
    "F731" = 1
    "G731" = 514.54
    "B732" = 65362
    "F732" = 76
    "G732" = 3106.12
    "B733" = 65079
    "F733" = 21
    "G733" = 858.27
    "F734" = 9
    "G734" = 424.89
    "F741" = 160
    "G741" = 38688
    "F742" = 142
    "G742" = 8092.58
    "B743" = 131525.59
    "F768" = 3763
    "G768" = 613782.9300000001
    "F771" = 545
    "G771" = 78834.25
    "B775" = 926160.51
    "B793" = 3881927.51
    "B794" = 3881927.51
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
